# Import Sales Lines from Excel
# The "Unit Cost" for the PARIS Guest Chair row (E3) was a staging value
# left over from the import; clear it back out and leave the selection
# where the user's cursor ended up after reviewing that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

$ws.Range("E3").ClearContents()
$ws.Range("E3").Select()
